$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-89 down to 42-90.
# Excel's row Insert() automatically carries down the formatting (e.g. the
# date number format s="2" on column D) from the surrounding rows.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with its values.
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44494
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112012
$ws.Range("G41").Value = "Espinaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 20
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = 8000
$ws.Range("N41").Value = "$/docena de atados"
$ws.Range("O41").Value = "Región de La Araucanía"
$ws.Range("P41").Value = 2667
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = "Hortaliza"
